$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B403").Value = 432885.0697202097
$ws.Range("B404").Value = 446023.4755036082
$ws.Range("B405").Value = 459212.0194237343
$ws.Range("B406").Value = 472448.2850765451
$ws.Range("B407").Value = 485729.8670802101
$ws.Range("B408").Value = 499054.375568775
$ws.Range("B409").Value = 512419.4406915503
$ws.Range("B410").Value = 525822.7171072613
$ws.Range("B411").Value = 539261.8884618283
$ws.Range("B412").Value = 552734.6718383151
$ws.Range("B413").Value = 565562.9827916964
$ws.Range("B414").Value = 578417.459756285
$ws.Range("B415").Value = 591295.9862221478
$ws.Range("B416").Value = 604196.4921922461
$ws.Range("B417").Value = 617116.9580977663
$ws.Range("B418").Value = 630055.4186166269
$ws.Range("B419").Value = 643009.9663813858
$ws.Range("B420").Value = 655978.7555620604
$ws.Range("B421").Value = 668960.0053085858
$ws.Range("B422").Value = 681952.0030367494
$ws.Range("B423").Value = 694953.1075405178
$ws.Range("B424").Value = 707961.7519126042
$ws.Range("B425").Value = 720976.4462540812
$ws.Range("B426").Value = 733995.7801526334
$ws.Range("B427").Value = 747018.4249078848
$ws.Range("B428").Value = 760043.1354809916
$ws.Range("B429").Value = 773068.7521444299
$ws.Range("B430").Value = 786094.2018067047
$ws.Range("B431").Value = 799118.4989854943
$ws.Range("B432").Value = 812140.7464016241
$ws.Range("B433").Value = 825160.1351652655
$ws.Range("B434").Value = 838175.9445248597
$ws.Range("B435").Value = 851187.5411486166
$ws.Range("B436").Value = 864194.3779079535
$ws.Range("B437").Value = 877195.9921321145
$ws.Range("B438").Value = 890192.0033033415
$ws.Range("B439").Value = 903182.1101625571
$ws.Range("B440").Value = 916166.0871964858
$ws.Range("B441").Value = 929143.7804786447
$ws.Range("B442").Value = 942115.1028386605
$ws.Range("B443").Value = 955080.0283369797
$ws.Range("B444").Value = 968038.5860253117
$ws.Range("B445").Value = 980990.8529770701
$ws.Range("B446").Value = 993936.9465767115
$ws.Range("B447").Value = 1006877.016062242
$ws.Range("B448").Value = 1019811.233321269
$ws.Range("B449").Value = 1032739.782947849
$ws.Range("B450").Value = 1045662.851574917
$ws.Range("B451").Value = 1058580.616505422
$ws.Range("B452").Value = 1071493.233674219
